$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '245.66'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = '22'
$ws.Cells.Item(2,7).Style = "Normal"
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '23.90'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = '22'
$ws.Cells.Item(3,7).Style = "Normal"
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '5.352'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = '22'
$ws.Cells.Item(4,7).Style = "Normal"
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '0.05822'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = '22'
$ws.Cells.Item(5,7).Style = "Normal"
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '3.371'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = '22'
$ws.Cells.Item(6,7).Style = "Normal"
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '6.466'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = '22'
$ws.Cells.Item(7,7).Style = "Normal"
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.8097'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = '22'
$ws.Cells.Item(8,7).Style = "Normal"
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.9192'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = '22'
$ws.Cells.Item(9,7).Style = "Normal"
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.1404'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = '22'
$ws.Cells.Item(10,7).Style = "Normal"
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.07380'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = '22'
$ws.Cells.Item(11,7).Style = "Normal"
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.03162'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,7).NumberFormat = "@"
$ws.Cells.Item(12,7).Value = '22'
$ws.Cells.Item(12,7).Style = "Normal"
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '0.03074'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,7).NumberFormat = "@"
$ws.Cells.Item(13,7).Value = '22'
$ws.Cells.Item(13,7).Style = "Normal"
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '0.09379'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,7).NumberFormat = "@"
$ws.Cells.Item(14,7).Value = '22'
$ws.Cells.Item(14,7).Style = "Normal"
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '3.846'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,7).NumberFormat = "@"
$ws.Cells.Item(15,7).Value = '22'
$ws.Cells.Item(15,7).Style = "Normal"
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '0.001556'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,7).NumberFormat = "@"
$ws.Cells.Item(16,7).Value = '22'
$ws.Cells.Item(16,7).Style = "Normal"
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '0.04680'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,7).NumberFormat = "@"
$ws.Cells.Item(17,7).Value = '22'
$ws.Cells.Item(17,7).Style = "Normal"
$ws.Cells.Item(18,2).NumberFormat = "@"
$ws.Cells.Item(18,2).Value = 'TigerCash'
$ws.Cells.Item(18,2).Style = "Normal"
$ws.Cells.Item(18,3).NumberFormat = "@"
$ws.Cells.Item(18,3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(18,3).Style = "Normal"
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '0.006045'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = '17TigerCashTCH'
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(18,7).NumberFormat = "@"
$ws.Cells.Item(18,7).Value = '22'
$ws.Cells.Item(18,7).Style = "Normal"
$ws.Cells.Item(19,2).NumberFormat = "@"
$ws.Cells.Item(19,2).Value = 'BitKan'
$ws.Cells.Item(19,2).Style = "Normal"
$ws.Cells.Item(19,3).NumberFormat = "@"
$ws.Cells.Item(19,3).Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Cells.Item(19,3).Style = "Normal"
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '0.001240'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = '18BitKanKAN'
$ws.Cells.Item(19,5).Style = "Normal"
$ws.Cells.Item(19,7).NumberFormat = "@"
$ws.Cells.Item(19,7).Value = '22'
$ws.Cells.Item(19,7).Style = "Normal"
$ws.Cells.Item(20,2).NumberFormat = "@"
$ws.Cells.Item(20,2).Value = 'HotbitToken'
$ws.Cells.Item(20,2).Style = "Normal"
$ws.Cells.Item(20,3).NumberFormat = "@"
$ws.Cells.Item(20,3).Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Cells.Item(20,3).Style = "Normal"
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '0.004689'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = '19HotbitTokenHTB'
$ws.Cells.Item(20,5).Style = "Normal"
$ws.Cells.Item(20,7).NumberFormat = "@"
$ws.Cells.Item(20,7).Value = '22'
$ws.Cells.Item(20,7).Style = "Normal"
$ws.Cells.Item(21,2).NumberFormat = "@"
$ws.Cells.Item(21,2).Value = 'NitroEx'
$ws.Cells.Item(21,2).Style = "Normal"
$ws.Cells.Item(21,3).NumberFormat = "@"
$ws.Cells.Item(21,3).Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Cells.Item(21,3).Style = "Normal"
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '0.00008796'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = '20NitroExNTX'
$ws.Cells.Item(21,5).Style = "Normal"
$ws.Cells.Item(21,7).NumberFormat = "@"
$ws.Cells.Item(21,7).Value = '22'
$ws.Cells.Item(21,7).Style = "Normal"
$ws.Cells.Item(22,2).NumberFormat = "@"
$ws.Cells.Item(22,2).Value = 'LEO'
$ws.Cells.Item(22,2).Style = "Normal"
$ws.Cells.Item(22,3).NumberFormat = "@"
$ws.Cells.Item(22,3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(22,3).Style = "Normal"
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '3.596'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = '21LEOLEO'
$ws.Cells.Item(22,5).Style = "Normal"
$ws.Cells.Item(22,7).NumberFormat = "@"
$ws.Cells.Item(22,7).Value = '22'
$ws.Cells.Item(22,7).Style = "Normal"
$ws.Cells.Item(23,2).NumberFormat = "@"
$ws.Cells.Item(23,2).Value = 'BTSEToken'
$ws.Cells.Item(23,2).Style = "Normal"
$ws.Cells.Item(23,3).NumberFormat = "@"
$ws.Cells.Item(23,3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(23,3).Style = "Normal"
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '2.150'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = '22BTSETokenBTSE'
$ws.Cells.Item(23,5).Style = "Normal"
$ws.Cells.Item(23,7).NumberFormat = "@"
$ws.Cells.Item(23,7).Value = '22'
$ws.Cells.Item(23,7).Style = "Normal"
$ws.Cells.Item(24,2).NumberFormat = "@"
$ws.Cells.Item(24,2).Value = 'One'
$ws.Cells.Item(24,2).Style = "Normal"
$ws.Cells.Item(24,3).NumberFormat = "@"
$ws.Cells.Item(24,3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(24,3).Style = "Normal"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '0.01069'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = '23OneONEBestin24h'
$ws.Cells.Item(24,5).Style = "Normal"
$ws.Cells.Item(24,7).NumberFormat = "@"
$ws.Cells.Item(24,7).Value = '22'
$ws.Cells.Item(24,7).Style = "Normal"
$ws.Cells.Item(25,7).NumberFormat = "@"
$ws.Cells.Item(25,7).Value = '22'
$ws.Cells.Item(25,7).Style = "Normal"
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '0.1319'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,7).NumberFormat = "@"
$ws.Cells.Item(26,7).Value = '22'
$ws.Cells.Item(26,7).Style = "Normal"
$ws.Cells.Item(27,7).NumberFormat = "@"
$ws.Cells.Item(27,7).Value = '22'
$ws.Cells.Item(27,7).Style = "Normal"
$ws.Cells.Item(28,7).NumberFormat = "@"
$ws.Cells.Item(28,7).Value = '22'
$ws.Cells.Item(28,7).Style = "Normal"
$ws.Cells.Item(29,7).NumberFormat = "@"
$ws.Cells.Item(29,7).Value = '22'
$ws.Cells.Item(29,7).Style = "Normal"
$ws.Cells.Item(30,7).NumberFormat = "@"
$ws.Cells.Item(30,7).Value = '22'
$ws.Cells.Item(30,7).Style = "Normal"
$ws.Cells.Item(31,7).NumberFormat = "@"
$ws.Cells.Item(31,7).Value = '22'
$ws.Cells.Item(31,7).Style = "Normal"
$ws.Cells.Item(32,7).NumberFormat = "@"
$ws.Cells.Item(32,7).Value = '22'
$ws.Cells.Item(32,7).Style = "Normal"
$ws.Cells.Item(33,7).NumberFormat = "@"
$ws.Cells.Item(33,7).Value = '22'
$ws.Cells.Item(33,7).Style = "Normal"
$ws.Cells.Item(34,7).NumberFormat = "@"
$ws.Cells.Item(34,7).Value = '22'
$ws.Cells.Item(34,7).Style = "Normal"
$ws.Cells.Item(35,7).NumberFormat = "@"
$ws.Cells.Item(35,7).Value = '22'
$ws.Cells.Item(35,7).Style = "Normal"
$ws.Cells.Item(36,7).NumberFormat = "@"
$ws.Cells.Item(36,7).Value = '22'
$ws.Cells.Item(36,7).Style = "Normal"
$ws.Cells.Item(37,7).NumberFormat = "@"
$ws.Cells.Item(37,7).Value = '22'
$ws.Cells.Item(37,7).Style = "Normal"
$ws.Cells.Item(38,7).NumberFormat = "@"
$ws.Cells.Item(38,7).Value = '22'
$ws.Cells.Item(38,7).Style = "Normal"
$ws.Cells.Item(39,7).NumberFormat = "@"
$ws.Cells.Item(39,7).Value = '22'
$ws.Cells.Item(39,7).Style = "Normal"
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.03847'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,7).NumberFormat = "@"
$ws.Cells.Item(40,7).Value = '22'
$ws.Cells.Item(40,7).Style = "Normal"
$ws.Cells.Item(41,2).NumberFormat = "@"
$ws.Cells.Item(41,2).Value = 'KickToken'
$ws.Cells.Item(41,2).Style = "Normal"
$ws.Cells.Item(41,3).NumberFormat = "@"
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Cells.Item(41,3).Style = "Normal"
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.006387'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = '40KickTokenKICK'
$ws.Cells.Item(41,5).Style = "Normal"
$ws.Cells.Item(41,7).NumberFormat = "@"
$ws.Cells.Item(41,7).Value = '22'
$ws.Cells.Item(41,7).Style = "Normal"
$ws.Cells.Item(42,2).NumberFormat = "@"
$ws.Cells.Item(42,2).Value = 'BKEXToken'
$ws.Cells.Item(42,2).Style = "Normal"
$ws.Cells.Item(42,3).NumberFormat = "@"
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Cells.Item(42,3).Style = "Normal"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.1066'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = '41BKEXTokenBKK'
$ws.Cells.Item(42,5).Style = "Normal"
$ws.Cells.Item(42,7).NumberFormat = "@"
$ws.Cells.Item(42,7).Value = '22'
$ws.Cells.Item(42,7).Style = "Normal"
$ws.Cells.Item(43,2).NumberFormat = "@"
$ws.Cells.Item(43,2).Value = 'CEJI'
$ws.Cells.Item(43,2).Style = "Normal"
$ws.Cells.Item(43,3).NumberFormat = "@"
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Cells.Item(43,3).Style = "Normal"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.003098'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = '42CEJICEJI'
$ws.Cells.Item(43,5).Style = "Normal"
$ws.Cells.Item(43,7).NumberFormat = "@"
$ws.Cells.Item(43,7).Value = '22'
$ws.Cells.Item(43,7).Style = "Normal"
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '0.009041'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,7).NumberFormat = "@"
$ws.Cells.Item(44,7).Value = '22'
$ws.Cells.Item(44,7).Style = "Normal"
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.00005246'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,7).NumberFormat = "@"
$ws.Cells.Item(45,7).Value = '22'
$ws.Cells.Item(45,7).Style = "Normal"
$ws.Cells.Item(46,7).NumberFormat = "@"
$ws.Cells.Item(46,7).Value = '22'
$ws.Cells.Item(46,7).Style = "Normal"
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.6852'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,7).NumberFormat = "@"
$ws.Cells.Item(47,7).Value = '22'
$ws.Cells.Item(47,7).Style = "Normal"
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.001825'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = '47BOLOBOLOWorstin24h'
$ws.Cells.Item(48,5).Style = "Normal"
$ws.Cells.Item(48,7).NumberFormat = "@"
$ws.Cells.Item(48,7).Value = '22'
$ws.Cells.Item(48,7).Style = "Normal"
$ws.Cells.Item(49,7).NumberFormat = "@"
$ws.Cells.Item(49,7).Value = '22'
$ws.Cells.Item(49,7).Style = "Normal"
$ws.Cells.Item(50,7).NumberFormat = "@"
$ws.Cells.Item(50,7).Value = '22'
$ws.Cells.Item(50,7).Style = "Normal"
$ws.Cells.Item(51,7).NumberFormat = "@"
$ws.Cells.Item(51,7).Value = '22'
$ws.Cells.Item(51,7).Style = "Normal"
